# "røret er vendt og data is gemt" - the rig was flipped around and a
# second run's data was saved: rename the original sheet to "Forsøg 1",
# add a new "Forsøg 2" sheet with a second experiment's readings, and
# relabel / extend the shared header strings with units + the new
# Delta_T column.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: rename the existing sheet, keep its data, add Delta_T ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Forsøg 1"

# --- Sheet 2: brand-new sheet placed right after "Forsøg 1" ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Forsøg 2"

# --- Header row (shared strings get new wording + a 4th header) ---
$ws1.Range("A1").Value = "M"
$ws1.Range("B1").Value = "T_start [C]"
$ws1.Range("C1").Value = "T_slut [C]"
$ws1.Range("D1").Value = "Delta_T [K]"

$ws2.Range("A1").Value = "M"
$ws2.Range("B1").Value = "T_start [C]"
$ws2.Range("C1").Value = "T_slut [C]"
$ws2.Range("D1").Value = "Delta_T [K]"

# --- Forsøg 1 measurements (rows 2-8) ---
$sheet1Data = @(
    @(100, 26.9, 28),
    @(100, 28.3, 29),
    @(100, 28.4, 29.6),
    @(50,  28.5, 29.2),
    @(52,  28.7, 29.2),
    @(50,  26.5, 27.6),
    @(26,  25.6, 26.9)
)

$row = 2
foreach ($r in $sheet1Data) {
    $ws1.Cells.Item($row, 1).Value = $r[0]
    $ws1.Cells.Item($row, 2).Value = $r[1]
    $ws1.Cells.Item($row, 3).Value = $r[2]
    $ws1.Range("D$row").Formula = "=C$row-B$row"
    $row++
}

# Rows 9-20: only the Delta_T formula was filled down (no raw readings)
for ($row = 9; $row -le 20; $row++) {
    $ws1.Range("D$row").Formula = "=C$row-B$row"
}

# --- Forsøg 2 measurements (rows 2-13) ---
$sheet2Data = @(
    @(100, 25.2, 28.4),
    @(100, 25.7, 28.4),
    @(100, 25.7, 28.2),
    @(150, 25.7, 29.4),
    @(150, 25.2, 29.4),
    @(150, 25.2, 29.2),
    @(50,  24.7, 26.7),
    @(50,  25.2, 26.8),
    @(50,  25.2, 26.3),
    @(20,  25,   25.4),
    @(20,  25.2, 25.9),
    @(20,  25.4, 25.9)
)

$row = 2
foreach ($r in $sheet2Data) {
    $ws2.Cells.Item($row, 1).Value = $r[0]
    $ws2.Cells.Item($row, 2).Value = $r[1]
    $ws2.Cells.Item($row, 3).Value = $r[2]
    $ws2.Range("D$row").Formula = "=C$row-B$row"
    $row++
}

# --- Selections / active sheet, matching the saved view state ---
$ws1.Range("G13").Select() | Out-Null
$ws2.Range("F5").Select() | Out-Null
$ws2.Activate() | Out-Null
